# Update cryptos list with latest scraped values (GitHub Actions refresh).
# Cells in column D that look like plain numbers are prefixed with a literal
# leading apostrophe so Excel stores them as text (matching the source feed's
# exact formatting, e.g. trailing zeros like "214.60") instead of coercing them
# into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.069.18'
$ws.Range('E2').Value = '  +0.35%  '
# Row 3
$ws.Range('D3').Value = '1.639.45'
$ws.Range('E3').Value = '  -0.02%  '
# Row 4
$ws.Range('E4').Value = '  +0.38%  '
# Row 5
$ws.Range('D5').Value = '''214.60'
$ws.Range('E5').Value = '  -0.56%  '
# Row 6
$ws.Range('D6').Value = '''0.504'
$ws.Range('E6').Value = '  -0.36%  '
# Row 7
$ws.Range('E7').Value = '  +0.34%  '
# Row 8
$ws.Range('D8').Value = '''0.251'
$ws.Range('E8').Value = '  -2.05%  '
# Row 9
$ws.Range('D9').Value = '''0.0625'
$ws.Range('E9').Value = '  -1.98%  '
# Row 10
$ws.Range('D10').Value = '''18.61'
$ws.Range('E10').Value = '  -5.02%  '
# Row 11
$ws.Range('E11').Value = '  -0.04%  '
# Row 12
$ws.Range('D12').Value = '1.747.75'
$ws.Range('E12').Value = '  +6.35%  '
# Row 13
$ws.Range('E13').Value = '  -1.68%  '
# Row 14
$ws.Range('D14').Value = '''0.531'
$ws.Range('E14').Value = '  -2.39%  '
# Row 15
$ws.Range('D15').Value = '''62.38'
$ws.Range('E15').Value = '  -1.02%  '
# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0₃0748'
$ws.Range('E16').Value = '  -2.00%  '
# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '26.073.47'
$ws.Range('E17').Value = '  +0.35%  '
# Row 18
$ws.Range('E18').Value = '  +0.36%  '
# Row 19
$ws.Range('D19').Value = '''190.77'
$ws.Range('E19').Value = '  -1.09%  '
# Row 20
$ws.Range('D20').Value = '''4.28'
$ws.Range('E20').Value = '  -1.89%  '
# Row 21
$ws.Range('D21').Value = '''9.59'
$ws.Range('E21').Value = '  -3.38%  '
# Row 22
$ws.Range('D22').Value = '''6.15'
$ws.Range('E22').Value = '  -2.04%  '
# Row 23
$ws.Range('D23').Value = '''144.24'
$ws.Range('E23').Value = '  +0.48%  '
# Row 24
$ws.Range('E24').Value = '  -0.95%  '
# Row 25
$ws.Range('E25').Value = '  +0.29%  '
# Row 26
$ws.Range('E26').Value = '  -1.40%  '
# Row 27
$ws.Range('E27').Value = '  -1.81%  '
# Row 28
$ws.Range('E28').Value = '  -2.30%  '
# Row 30
$ws.Range('E30').Value = '  -3.35%  '
# Row 31
$ws.Range('E31').Value = '  -2.20%  '
# Row 32
$ws.Range('E32').Value = '  -3.57%  '
# Row 33
$ws.Range('E33').Value = '  -0.48%  '
# Row 34
$ws.Range('E34').Value = '  -1.51%  '
# Row 35
$ws.Range('D35').Value = '''0.878'
$ws.Range('E35').Value = '  -2.48%  '
# Row 36
$ws.Range('D36').Value = '1.122.24'
$ws.Range('E36').Value = '  -1.01%  '
# Row 37
$ws.Range('E37').Value = '  -0.19%  '
# Row 38
$ws.Range('D38').Value = '''0.524'
$ws.Range('E38').Value = '  -3.57%  '
# Row 39
$ws.Range('E39').Value = '  -1.51%  '
# Row 40
$ws.Range('D40').Value = '''98.77'
$ws.Range('E40').Value = '  -0.50%  '
# Row 41
$ws.Range('E41').Value = '  -1.61%  '
# Row 42
$ws.Range('E42').Value = '  -3.43%  '
# Row 43
$ws.Range('E43').Value = '  -0.63%  '
# Row 44
$ws.Range('D44').Value = '''55.22'
$ws.Range('E44').Value = '  -2.57%  '
# Row 45
$ws.Range('D45').Value = '''0.0523'
$ws.Range('E45').Value = '  -1.30%  '
# Row 46
$ws.Range('D46').Value = '''1.48'
$ws.Range('E46').Value = '  +0.71%  '
# Row 47
$ws.Range('E47').Value = '  -0.10%  '
# Row 48
$ws.Range('D48').Value = '''7.59'
$ws.Range('E48').Value = '  -0.83%  '
# Row 49
$ws.Range('E49').Value = '  +0.12%  '
# Row 50
$ws.Range('D50').Value = '''0.0928'
$ws.Range('E50').Value = '  -3.44%  '
# Row 51
$ws.Range('E51').Value = '  -0.70%  '
